# "Estabilizacion pago tc propias"
# Replace the 3rd test row's own-card ("Personal Visa") payment data with a
# fresh own-card ("Personal American Express") payment scenario: new card
# alias, new user id, new payment value and a new destination account.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D4: numeroDocumento/usuario field -> "pagotdc1"
# Keep the destination cell's quote-prefixed "text" entry style (like its
# neighbour E4) rather than the original bold/quote-prefixed style.
$ws.Range("E4").Copy()
$ws.Range("D4").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("D4").Formula = "'pagotdc1"

# O4/P4/R4: switch the card from "Personal Visa" / *5880 / 480369
# to "Personal American Express" / *0702 / 1000
$ws.Range("O4").Value = "Personal American Express"
$ws.Range("P4").Value = "*0702"
$ws.Range("R4").Value = "1000"

# U4: destination account number -> 406-738430-25
$ws.Range("U4").Value = "406-738430-25"
$ws.Range("U4").Font.Name = "Calibri"

$ws.Range("U4").Select() | Out-Null
